$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.167.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.286.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.73%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.73"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -9.58%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.289.46"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.466"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.72%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.847.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.80%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -8.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.292.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.080.31"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "368.67"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.09%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.38"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.526"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.431.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0983"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -11.79%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.171"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.06%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -8.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.11"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.26"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.18%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.76"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.47%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.64%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.52"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.83%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.314.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.96%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0711"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.51%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.69"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.37%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.94"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -18.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.741"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.81%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.05"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.03%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.297.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.29"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.10%  "
